$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-23 Tuesday" "2024-07-24 Wednesday"

Replace-Text "126÷7=" "956÷3="
Replace-Text "540÷7=" "911÷9="
Replace-Text "737÷6=" "470÷8="
Replace-Text "450÷9=" "211÷8="
Replace-Text "788÷4=" "446÷6="

Replace-Text "751÷3=" "874÷5="
Replace-Text "977÷9=" "960÷8="
Replace-Text "244÷5=" "359÷4="
Replace-Text "897÷7=" "676÷5="
Replace-Text "552÷8=" "151÷9="

Replace-Text "904÷5=" "268÷5="
Replace-Text "755÷2=" "913÷9="
Replace-Text "608÷9=" "959÷8="
Replace-Text "728÷7=" "964÷5="
Replace-Text "114÷6=" "288÷9="

Replace-Text "859÷5=" "445÷7="
Replace-Text "251÷4=" "929÷9="
Replace-Text "858÷9=" "772÷6="
Replace-Text "627÷5=" "234÷4="
Replace-Text "846÷2=" "273÷4="

Replace-Text "958÷4=" "649÷5="
Replace-Text "661÷7=" "302÷3="
Replace-Text "568÷9=" "225÷7="
Replace-Text "966÷4=" "352÷8="
Replace-Text "107÷3=" "256÷6="
